$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking / percent-looking text values that must
# remain plain text (as in the source data), so force text format before writing,
# then restore the default (unstyled) cell style so no visual formatting changes.
$deRange = $ws.Range("D2:E51")
$deRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '311.69'
$ws.Cells.Item(2, 5).Value = '-2.88%'
$ws.Cells.Item(3, 4).Value = '53.86'
$ws.Cells.Item(3, 5).Value = '9.33%'
$ws.Cells.Item(4, 4).Value = '5.103'
$ws.Cells.Item(4, 5).Value = '-4.26%'
$ws.Cells.Item(5, 4).Value = '0.07900'
$ws.Cells.Item(5, 5).Value = '-1.92%'
$ws.Cells.Item(6, 4).Value = '4.569'
$ws.Cells.Item(6, 5).Value = '-0.89%'
$ws.Cells.Item(7, 4).Value = '1.395'
$ws.Cells.Item(7, 5).Value = '4.24%'
$ws.Cells.Item(8, 4).Value = '1.660'
$ws.Cells.Item(8, 5).Value = '0.85%'
$ws.Cells.Item(9, 4).Value = '0.1236'
$ws.Cells.Item(9, 5).Value = '-2.09%'
$ws.Cells.Item(10, 4).Value = '0.2001'
$ws.Cells.Item(10, 5).Value = '1.56%'
$ws.Cells.Item(11, 4).Value = '0.04726'
$ws.Cells.Item(11, 5).Value = '0.20%'
$ws.Cells.Item(12, 4).Value = '0.09431'
$ws.Cells.Item(12, 5).Value = '-2.47%'
$ws.Cells.Item(13, 4).Value = '0.1043'
$ws.Cells.Item(13, 5).Value = '-0.41%'
$ws.Cells.Item(14, 4).Value = '0.001278'
$ws.Cells.Item(14, 5).Value = '-3.97%'
$ws.Cells.Item(15, 4).Value = '0.005817'
$ws.Cells.Item(15, 5).Value = '-1.25%'
$ws.Cells.Item(16, 4).Value = '3.334'
$ws.Cells.Item(16, 5).Value = '-0.27%'
$ws.Cells.Item(17, 4).Value = '2.436'
$ws.Cells.Item(17, 5).Value = '0.05%'
$ws.Cells.Item(18, 4).Value = '0.3480'
$ws.Cells.Item(18, 5).Value = '-1.20%'
$ws.Cells.Item(19, 4).Value = '8.369'
$ws.Cells.Item(19, 5).Value = '4.23%'
$ws.Cells.Item(20, 4).Value = '0.1358'
$ws.Cells.Item(20, 5).Value = '-1.64%'
$ws.Cells.Item(21, 4).Value = '0.2918'
$ws.Cells.Item(21, 5).Value = '-5.74%'
$ws.Cells.Item(22, 4).Value = '0.04179'
$ws.Cells.Item(22, 5).Value = '-0.74%'
$ws.Cells.Item(23, 4).Value = '0.001256'
$ws.Cells.Item(23, 5).Value = '-3.62%'
$ws.Cells.Item(24, 4).Value = '0.003992'
$ws.Cells.Item(24, 5).Value = '-7.48%'
$ws.Cells.Item(25, 4).Value = '0.0001346'
$ws.Cells.Item(25, 5).Value = '-0.35%'
$ws.Cells.Item(26, 4).Value = '0.0003530'
$ws.Cells.Item(26, 5).Value = '-0.30%'
$ws.Cells.Item(38, 4).Value = '0.02639'
$ws.Cells.Item(38, 5).Value = '-2.61%'
$ws.Cells.Item(39, 4).Value = '0.05933'
$ws.Cells.Item(39, 5).Value = '-0.51%'
$ws.Cells.Item(40, 5).Value = '0.16%'
$ws.Cells.Item(41, 4).Value = '0.1695'
$ws.Cells.Item(41, 5).Value = '15.79%'
$ws.Cells.Item(42, 4).Value = '0.007963'
$ws.Cells.Item(42, 5).Value = '-0.72%'
$ws.Cells.Item(43, 4).Value = '0.008212'
$ws.Cells.Item(43, 5).Value = '3.71%'
$ws.Cells.Item(44, 4).Value = '0.008349'
$ws.Cells.Item(44, 5).Value = '5.83%'
$ws.Cells.Item(45, 4).Value = '0.3445'
$ws.Cells.Item(45, 5).Value = '-1.59%'
$ws.Cells.Item(46, 4).Value = '0.00007224'
$ws.Cells.Item(46, 5).Value = '4.60%'
$ws.Cells.Item(47, 5).Value = '-0.32%'
$ws.Cells.Item(48, 4).Value = '0.002613'
$ws.Cells.Item(48, 5).Value = '-34.70%'
$ws.Cells.Item(49, 4).Value = '0.05550'
$ws.Cells.Item(49, 5).Value = '0.63%'
$ws.Cells.Item(50, 4).Value = '0.00002094'
$ws.Cells.Item(50, 5).Value = '-0.32%'
$ws.Cells.Item(51, 5).Value = '-0.32%'

$deRange.Style = "Normal"

# Coin name / link swaps (plain text columns)
$ws.Cells.Item(11, 2).Value = 'BitrueCoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(48, 2).Value = 'CoinbaseStockToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Cells.Item(49, 2).Value = 'BOLO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
